## Documentation update: add a note about flags=global-headers for ffenc_mpeg4
## and insert the missing WMV8 encoder table row, plus note WMV extension.

$d = $word.ActiveDocument

## --- Table 1: "List of supported VIDEO ENCODER types" ---
## Current (before) rows, relevant slice:
##   8  MPEG-4 part 2           | ffenc_mpeg4      | bitrate | Bitrate in bit/s. Default is 300 000.
##   9  MPEG-2 video            | ffenc_mpeg2video | bitrate | Bitrate in bit/s. Default is 300 000.
##   10 Windows Media Video 8   | ffenc_wmv2       | bitrate | Bitrate in bit/s. Default is 300 000.
##   11 Flash video             | ffenc_flv        | bitrate | Bitrate in bit/s. Default is 300 000.
##
## Target (after):
##   8  MPEG-4 part 2           | ffenc_mpeg4      | flags   | Specify ‘flags=global-headers’ when using with .3gp or .mp4 formats.
##   9  (row 9 Name/Encoder cells emptied - becomes a second options row for MPEG-4 part 2)
##   10 MPEG-2 video            | ffenc_mpeg2video | bitrate | Bitrate in bit/s. Default is 300 000.
##   11 Windows Media Video 8   | ffenc_wmv2       | bitrate | Bitrate in bit/s. Default is 300 000.  (new row)
##   12 Flash video             | ffenc_flv        | bitrate | Bitrate in bit/s. Default is 300 000.

$t1 = $d.Tables.Item(1)

## 1) Insert a brand-new row right before the "Flash video" row (currently row 11).
##    It will hold what used to be the "Windows Media Video 8" row's content.
$t1.Rows.Add($t1.Rows.Item(11)) | Out-Null

## 2) Populate the freshly inserted row (now row 11) with the WMV8 encoder info
##    (this is the content that used to live in row 10 before the edits below).
$t1.Cell(11, 1).Range.Text = "Windows Media Video 8"
$t1.Cell(11, 2).Range.Text = "ffenc_wmv2"
$t1.Cell(11, 3).Range.Text = "bitrate"
$t1.Cell(11, 4).Range.Text = "Bitrate in bit/s. Default is 300 000."

## 3) Retarget the old "Windows Media Video 8" row (row 10) to describe MPEG-2 video instead.
$t1.Cell(10, 1).Range.Text = "MPEG-2 video"
$t1.Cell(10, 2).Range.Text = "ffenc_mpeg2video"
## columns 3/4 of row 10 ("bitrate" / "Bitrate in bit/s. Default is 300 000.") stay as-is.

## 4) Clear out the Name/Encoder cells of the old "MPEG-2 video" row (row 9) - this becomes
##    a second options row for "MPEG-4 part 2" above it (empty Name/Encoder cells).
$c91 = $t1.Cell(9, 1).Range
$d.Range($c91.Start, $c91.End - 1).Delete() | Out-Null
$c92 = $t1.Cell(9, 2).Range
$d.Range($c92.Start, $c92.End - 1).Delete() | Out-Null

## 5) Change the option/description of row 9 to document the new "flags" option.
$t1.Cell(9, 3).Range.Text = "flags"
$t1.Cell(9, 4).Range.Text = "Specify ‘flags=global-headers’ when using with .3gp or .mp4 formats."

## --- Table 2: "List of supported CONTAINER types" ---
## Append ", .WMV" to the file-extension cell of the ASF row.
## (Cell.Range.Text carries trailing paragraph-mark/cell-mark control chars,
## so trim those off before comparing.)
$t2 = $d.Tables.Item(2)
for ($r = 1; $r -le $t2.Rows.Count; $r++) {
    $cellText = $t2.Cell($r, 3).Range.Text.TrimEnd([char]13, [char]7)
    if ($cellText -eq ".ASF") {
        $t2.Cell($r, 3).Range.Text = ".ASF, .WMV"
        break
    }
}
